# Disaggregation of commodity Copper
#
# "Copper ores and concentrates" is renamed to "Copper" everywhere it is
# used (column C, row 4, on every yearly worksheet). A handful of
# worksheets also carry an updated D4 total (last-digit precision refresh
# that rode along with the relabeling in the source data).

$wb = $excel.ActiveWorkbook

$oldLabel = "Copper ores and concentrates"
$newLabel = "Copper"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("C4")
    if ($cell.Value2 -eq $oldLabel) {
        $cell.Value = $newLabel
    }
}

# Worksheet index (1-based, tab order) -> refreshed D4 total.
$updatedTotals = @(
    ,@(31, 59078.62275518187)
    ,@(42, 211682.5044181577)
    ,@(49, 725988.578645583)
    ,@(59, 1859543.682973828)
    ,@(66, 818774.3147223982)
    ,@(74, 791405.3239931302)
    ,@(75, 932849.6529590308)
)

foreach ($row in $updatedTotals) {
    $idx = $row[0]
    $val = $row[1]
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("D4").Value = $val
}
